# Applies the "version 4 5 6" commit:
#  1. Fixes the BE sheet's Wave 4 received-date typo (2020 -> 2021).
#  2. Appends BE Wave 5 and Wave 6 rows.
#  3. Adds a new "Group2" worksheet (Lithuania/Finland/Switzerland wave-1
#     tracking), positioned between "Group1" and "Sheet1".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1 & 2: BE sheet - fix the Wave 4 date and append Wave 5 / Wave 6 rows
# ---------------------------------------------------------------------
$be = $wb.Worksheets.Item("BE")

# Wave 4's date_recieved was mistyped as 2020-01-11; correct it to 2021-01-11.
$be.Cells.Item(13, 7).Value = (Get-Date -Year 2021 -Month 1 -Day 11 -Hour 0 -Minute 0 -Second 0)

# Wave 5 (row 14)
$be.Cells.Item(14, 1).Value = "be"
$be.Cells.Item(14, 2).Value = 4
$be.Cells.Item(14, 3).Value = 0
$be.Cells.Item(14, 4).Value = 13
$be.Cells.Item(14, 5).Value = "B"
$be.Cells.Item(14, 6).Value = 5
$be.Cells.Item(14, 7).Value = (Get-Date -Year 2021 -Month 1 -Day 20 -Hour 0 -Minute 0 -Second 0)
$be.Cells.Item(14, 8).Value = "20_060765_BE2_Wave5_Final_v1_200121_IntClientUse"
$be.Cells.Item(14, 9).Value = "be_wk12_20200111_pB_wv05"

# Wave 6 (row 15)
$be.Cells.Item(15, 1).Value = "be"
$be.Cells.Item(15, 2).Value = 4
$be.Cells.Item(15, 3).Value = 0
$be.Cells.Item(15, 4).Value = 13
$be.Cells.Item(15, 5).Value = "B"
$be.Cells.Item(15, 6).Value = 6
$be.Cells.Item(15, 7).Value = (Get-Date -Year 2021 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0)
$be.Cells.Item(15, 8).Value = "20_060765_BE2_Wave6_Final_v1_270121_IntClientUse"
$be.Cells.Item(15, 9).Value = "be_wk12_20200111_pB_wv06"

# Match the date format + style already used by the rest of column G.
$be.Range("G14:G15").NumberFormat = $be.Cells.Item(13, 7).NumberFormat

# ---------------------------------------------------------------------
# 3: Insert a new "Group2" worksheet just before "Sheet1"
# ---------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$group2 = $wb.Worksheets.Add($sheet1)
$group2.Name = "Group2"

$headers = @("country", "survey_version", "locked", "week", "panel", "wave", "date_recieved", "spss_name", "r_name", "r_saved")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $group2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2: Lithuania
$group2.Cells.Item(2, 1).Value = "lt"
$group2.Cells.Item(2, 2).Value = 6
$group2.Cells.Item(2, 3).Value = 0
$group2.Cells.Item(2, 4).Value = 1
$group2.Cells.Item(2, 5).Value = "A"
$group2.Cells.Item(2, 6).Value = 1
$group2.Cells.Item(2, 7).Value = (Get-Date -Year 2021 -Month 1 -Day 11 -Hour 0 -Minute 0 -Second 0)
$group2.Cells.Item(2, 8).Value = "20-030971_G2_Merged_Wave1_Final_v1_04022021_IntClientUse"
$group2.Cells.Item(2, 9).Formula = '=A2&"_"&"wk"&TEXT(D2,"00")&"_"&YEAR(G3)&TEXT(G3,"MM")&TEXT(G3,"DD")&"_p"&E2&"_wv"&TEXT(F2,"00")&""'

# Row 3: Finland
$group2.Cells.Item(3, 1).Value = "fi"
$group2.Cells.Item(3, 2).Value = 6
$group2.Cells.Item(3, 3).Value = 0
$group2.Cells.Item(3, 4).Value = 1
$group2.Cells.Item(3, 5).Value = "A"
$group2.Cells.Item(3, 6).Value = 1
$group2.Cells.Item(3, 7).Value = (Get-Date -Year 2021 -Month 1 -Day 12 -Hour 0 -Minute 0 -Second 0)
$group2.Cells.Item(3, 8).Value = "20-030971_G2_Merged_Wave1_Final_v1_04022021_IntClientUse"
$group2.Cells.Item(3, 9).Formula = '=A3&"_"&"wk"&TEXT(D3,"00")&"_"&YEAR(G4)&TEXT(G4,"MM")&TEXT(G4,"DD")&"_p"&E3&"_wv"&TEXT(F3,"00")&""'

# Row 4: Switzerland
$group2.Cells.Item(4, 1).Value = "ch"
$group2.Cells.Item(4, 2).Value = 6
$group2.Cells.Item(4, 3).Value = 0
$group2.Cells.Item(4, 4).Value = 1
$group2.Cells.Item(4, 5).Value = "A"
$group2.Cells.Item(4, 6).Value = 1
$group2.Cells.Item(4, 7).Value = (Get-Date -Year 2021 -Month 1 -Day 13 -Hour 0 -Minute 0 -Second 0)
$group2.Cells.Item(4, 8).Value = "20-030971_G2_Merged_Wave1_Final_v1_04022021_IntClientUse"
$group2.Cells.Item(4, 9).Formula = '=A4&"_"&"wk"&TEXT(D4,"00")&"_"&YEAR(G5)&TEXT(G5,"MM")&TEXT(G5,"DD")&"_p"&E4&"_wv"&TEXT(F4,"00")&""'

$group2.Range("G2:G4").NumberFormat = $be.Cells.Item(13, 7).NumberFormat

Write-Host "edit applied"
